{"js": "// Apply the tracked edits described by the diff:\n//  1. Remove the (stale) \"_GoBack\" bookmark that Word auto-manages around\n//     the last edited location.\n//  2. Normalize the \">>>  your stuff after this line >>>\" paragraph back\n//     into a single run (clearing the leftover grammar-check run splits /\n//     proofErr markers).\n//  3. Rename the signature line from \"Ben changing things up!\" to\n//     \"Tracy changing things up!\".\n\nconst body = context.document.body;\n\n// 1) Drop the _GoBack bookmark (start + end) if present.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Re-normalize the \">>>\" marker paragraph to a single, clean run.\nconst markerParas = body.search(\">>>  your stuff after this line >>>\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nmarkerParas.load(\"items\");\nawait context.sync();\n\nif (markerParas.items.length > 0) {\n  const markerParagraph = markerParas.items[0].paragraphs.getFirst();\n  markerParagraph.load(\"text\");\n  await context.sync();\n  const markerText = markerParagraph.text;\n  markerParagraph.clear();\n  await context.sync();\n  markerParagraph.insertText(markerText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Replace the author's name in the \"changing things up!\" line.\nconst nameResults = body.search(\"Ben changing things up!\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  const nameParagraph = nameResults.items[0].paragraphs.getFirst();\n  nameParagraph.clear();\n  await context.sync();\n  nameParagraph.insertText(\"Tracy changing things up!\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the tracked edits described by the diff:\n#  1. Remove the (stale) \"_GoBack\" bookmark that Word auto-manages around\n#     the last edited location.\n#  2. Normalize the \">>>  your stuff after this line >>>\" paragraph back\n#     into a single run (clearing the leftover grammar-check run splits /\n#     proofErr markers).\n#  3. Rename the signature line from \"Ben changing things up!\" to\n#     \"Tracy changing things up!\".\n\n$d = $word.ActiveDocument\n\n# 1) Drop the _GoBack bookmark (start + end) if present.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Re-normalize the \">>>\" marker paragraph to a single, clean run by\n#    finding it and \"replacing\" it with the same text -- this clears any\n#    leftover proofing-error run splits without touching the wording.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \">>>  your stuff after this line >>>\"\n$find.Replacement.Text = \">>>  your stuff after this line >>>\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 3) Replace the author's name in the \"changing things up!\" line.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Ben changing things up!\"\n$find2.Replacement.Text = \"Tracy changing things up!\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
